$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "66.947.95"
$ws.Range("E2").Value = "  +5.47%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.700.74"
$ws.Range("E3").Value = "  +6.66%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "424.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.690.77"
$ws.Range("E7").Value = "  +6.56%  "

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.640"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "

# Row 9: USDC
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10: Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.765"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.58%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.183"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.62%  "

# Row 12: ShibaInu
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000393"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +70.57%  "

# Row 13: Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.55%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.276.79"
$ws.Range("E15").Value = "  +5.95%  "

# Row 16: TRON
$ws.Range("E16").Value = "  +0.11%  "

# Row 17: Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.701.10"
$ws.Range("E18").Value = "  +7.04%  "

# Row 19: Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.78%  "

# Row 20: Polygon
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.74%  "

# Row 21: WrappedBTC
$ws.Range("D21").Value = "66.786.32"
$ws.Range("E21").Value = "  +5.31%  "

# Row 22: BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.12%  "

# Row 23: InternetComputer(DFINITY)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.28%  "

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.61%  "

# Row 25: ImmutableX
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.06%  "

# Row 26: EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.60%  "

# Row 27: Filecoin -> PancakeSwap
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.55%  "

# Row 28: PancakeSwap -> Filecoin
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29: LEO
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.78%  "

# Row 30: Toncoin -> Cosmos
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.77%  "

# Row 31: Cosmos -> Toncoin
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.62%  "

# Row 32: Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.77%  "

# Row 33: RenderToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.55%  "

# Row 34: InjectiveProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.90%  "

# Row 35: Kaspa
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.162"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.48%  "

# Row 36: Dai
$ws.Range("E36").Value = "  +0.09%  "

# Row 37: OKB
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.70%  "

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0489"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

# Row 39: ThetaToken -> PEPE
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0725"
$ws.Range("E39").Value = "  +13.39%  "

# Row 40: PEPE -> ThetaToken
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +29.52%  "

# Row 41: Stellar
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.146"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.23%  "

# Row 42: EnergySwap -> FirstDigitalUSD
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "

# Row 43: FirstDigitalUSD -> LidoDAOToken
$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.39%  "

# Row 44: LidoDAOToken -> EnergySwap
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +25.72%  "

# Row 45: Monero -> ARBITRUM
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.72%  "

# Row 46: WEMIXToken -> Stacks
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.93%  "

# Row 47: ARBITRUM -> Monero
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

# Row 48: Stacks -> NEARProtocol
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.71%  "

# Row 49: NEARProtocol -> WEMIXToken
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.55%  "

# Row 50: TheGraph
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.305"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.64%  "

# Row 51: Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.157"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.99%  "
